$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): Job posting column headers ---
$ws.Range("A1").Value = "Job_Id"
$ws.Range("B1").Value = "Job_Title"
$ws.Range("C1").Value = "Job_Description"
$ws.Range("D1").Value = "Total_Years_Min_Exp"
$ws.Range("E1").Value = "Total_Years_Max_Exp"
$ws.Range("F1").Value = "Linked_Posted"
$ws.Range("G1").Value = "Resume_received"
$ws.Range("H1").Value = "Resume_downloaded"

# Extend the existing header style (bold, bordered, centered - already used by
# A1:C1) across the newly added header cells D1:H1, by copying the format
# from an already-styled header cell.
$ws.Range("A1").Copy()
$ws.Range("D1:H1").PasteSpecial(-4122)  # xlPasteFormats

# --- Data row (row 2): first job posting record ---
$ws.Range("A2").Value = "JD_001"
$ws.Range("B2").Value = "Senior Engineer"
$ws.Range("C2").Value = "We are seeking a Software Engineer to build and maintain high-quality software solutions"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0

# A2 previously carried a wrap-text-only style left over from the blank
# template; clear it back to the plain/default style so the data row is
# unstyled like the rest of the record.
$ws.Range("Z100").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats

$excel.CutCopyMode = $false
